$wb = $excel.ActiveWorkbook

# Auto update stock data: bump the latest ("Date_1") row for each ticker
# sheet to 2025/10/28 and refresh the EBITDA figure that came with it.
$updates = @(
    @{ Sheet = "Alcoa";                       EBITDA = "4.85" },
    @{ Sheet = "Rio Tinto";                   EBITDA = "7.65" },
    @{ Sheet = "Norsk Hydro";                 EBITDA = $null },
    @{ Sheet = "Reliance Steel & Aluminum";   EBITDA = "12.33" },
    @{ Sheet = "Kaiser Aluminum";             EBITDA = "9.79" },
    @{ Sheet = "Ryerson Holding";             EBITDA = "20.35" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    $dateCell = $ws.Range("A2")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2025/10/28"

    if ($u.EBITDA -ne $null) {
        $ebitdaCell = $ws.Range("B2")
        $ebitdaCell.NumberFormat = "@"
        $ebitdaCell.Value = $u.EBITDA
    }
}
